$wb = $excel.ActiveWorkbook

# --- db_lines_out: zero out the two non-zero "month 12" cfa / prod-attained cells ---
$wsLines = $wb.Worksheets.Item("db_lines_out")
$wsLines.Range("L13").Value = 0
$wsLines.Range("N13").Value = 0
$wsLines.Range("L25").Value = 0
$wsLines.Range("N25").Value = 0

# --- db_sku_out: drop the trailing "month 12" duplicate rows for the first sku
#     (rows 24:25) and the trailing "month 12" block for the other skus
#     (rows 268:289). Delete bottom-up so earlier row numbers stay valid. ---
$wsSku = $wb.Worksheets.Item("db_sku_out")
$wsSku.Range("A268:A289").EntireRow.Delete()
$wsSku.Range("A24:A25").EntireRow.Delete()

# --- db_inventory_out: drop the trailing "month 12" block (rows 90:97) ---
$wsInv = $wb.Worksheets.Item("db_inventory_out")
$wsInv.Range("A90:A97").EntireRow.Delete()

# --- db_sku_logistics_out: drop the trailing "month 12" row (row 13) ---
$wsLog = $wb.Worksheets.Item("db_sku_logistics_out")
$wsLog.Range("A13").EntireRow.Delete()
